# Apply updated crypto price/volume data as per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.396.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.38%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.804.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").Value = "'  +0.13%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'306.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.55%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4518"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.46%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.46%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'46.25"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.20%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.07069"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.48%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.8886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.86%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07812"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.33%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'19.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.71%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'1.829.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.71%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'5.275"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.06%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'6.300"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.40%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'85.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.52%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  +0.07%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.000008499"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.45%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +0.07%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'26.452.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.22%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'14.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.39%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'4.960"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.068.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.46%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'10.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.28%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'1.960"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.78%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'152.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.52%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -0.35%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.066"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.82%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'112.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.72%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'4.843"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.65%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.08689"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.43%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'3.079"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.41%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'2.797"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +11.26%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'4.454"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.38%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -0.37%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -0.69%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.077"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.22%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.01933"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.22%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.912"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.47%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.05106"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.60%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.5064"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.24%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'6.774"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.69%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.1510"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.61%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'8.013"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.37%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.008"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.03%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.4661"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.42%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'10.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.39%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'100.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.16%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.575"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.15%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.05973"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.27%  "
$ws.Range("E51").Style = "Normal"
